# The due_date column (C) stores dates as plain text (e.g. "2024-01-01").
# Assigning a date-looking string straight to .Value would make Excel
# auto-convert it to a real date serial number, so we briefly force the
# cells to a text number format, assign the new text values, then restore
# the cells' style so the underlying formatting is unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("C2:C4")
$rng.NumberFormat = "@"

$ws.Range("C2").Value = "2024-01-15"
$ws.Range("C3").Value = "2024-02-15"
$ws.Range("C4").Value = "2024-03-15"

$rng.Style = "Normal"
